$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, shifting existing rows 25-49 down to 26-50
$ws.Rows(25).Insert()

# Populate the new row 25 with the new weekly record
$ws.Cells.Item(25, 1).Value = 9
$ws.Cells.Item(25, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(25, 3).Value = 'Metropolitana'
$ws.Cells.Item(25, 4).Value = 44763
$ws.Cells.Item(25, 5).Value = 13
$ws.Cells.Item(25, 6).Value = 100112035
$ws.Cells.Item(25, 7).Value = 'Bruselas (repollito)'
$ws.Cells.Item(25, 8).Value = 'Sin especificar'
$ws.Cells.Item(25, 9).Value = 'Primera'
$ws.Cells.Item(25, 10).Value = 34
$ws.Cells.Item(25, 11).Value = 20000
$ws.Cells.Item(25, 12).Value = 20000
$ws.Cells.Item(25, 13).Value = 20000
$ws.Cells.Item(25, 14).Value = '$/malla 15 kilos'
$ws.Cells.Item(25, 15).Value = 'Hijuelas'
$ws.Cells.Item(25, 16).Value = 1333
$ws.Cells.Item(25, 17).Value = 15
$ws.Cells.Item(25, 18).Value = 'Hortaliza'
